# refactor(Proveedores - Solicitudes - Órdenes Excel): Agregadas observaciones
#
# Adds a new "OBSERVACIÓN" column (column O) to the "Movimientos" sheet,
# mirroring the header formatting of the existing columns (A1:K1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Movimientos")

# New header cell in column O, row 1
$cell = $ws.Range("O1")
$cell.Value = "OBSERVACIÓN"

# Match the centered-header formatting used by the other header cells
# (style index "s=1" in the sheet, i.e. horizontally centered text)
$cell.HorizontalAlignment = -4108  # xlCenter

# Give the new column a sensible width (23 characters, as authored)
$ws.Columns.Item(15).ColumnWidth = 22.166666666666668
